$wb = $excel.ActiveWorkbook

$wsActors    = $wb.Worksheets.Item("Actors")
$wsRooms     = $wb.Worksheets.Item("Rooms")
$wsScripts   = $wb.Worksheets.Item("Scripts")
$wsVariables = $wb.Worksheets.Item("Variables")
$wsSound     = $wb.Worksheets.Item("Sound")
$wsObjects   = $wb.Worksheets.Item("Objects")

# --- Data edits -------------------------------------------------------
# New shared strings are appended in first-use order, so write the cells
# in the same order they appear in the target file so the shared string
# table indices line up: 768, 769, 770, 771.

# Objects!D7 - "08 = Not Available" (new string 768)
$wsObjects.Range("D7").Value = "08 = Not Available"

# Variables!G55 - "2 = Ed has package" (new string 769)
$wsVariables.Range("G55").Value = "2 = Ed has package"

# Insert new row 66 in Variables: id 87, "Weird Ed: Looking for plans" (new string 770)
$wsVariables.Rows.Item(66).Insert()
$wsVariables.Range("A66").Value = 87
$wsVariables.Range("C66").Value = "Weird Ed: Looking for plans"

# Insert new row 94 (post shift) in Variables: id 134, "Glass Jar has content" (new string 771)
$wsVariables.Rows.Item(94).Insert()
$wsVariables.Range("A94").Value = 134
$wsVariables.Range("C94").Value = "Glass Jar has content"

# --- View / selection state -------------------------------------------
$wsActors.Activate()
$wsActors.Range("B21").Select()

$wsRooms.Activate()
$wsRooms.Range("B39").Select()

$wsVariables.Activate()
$wsVariables.Range("E94").Select()

$wsSound.Activate()
$wsSound.Range("B32").Select()

$wsObjects.Activate()
$wsObjects.Range("C386").Select()

# Scripts must end up as the active tab (matches original activeTab=2)
$wsScripts.Activate()
$wsScripts.Range("D100").Select()
